$d = $word.ActiveDocument

# Locate the target paragraph: the empty "List Bullet" paragraph that
# immediately follows the "Space Complexity: O(n)" line (between the
# SortArrayByParity and TwoSum sections). It currently carries manual
# bold formatting on its paragraph mark and no explicit indent.
$target = $null
foreach ($para in $d.Paragraphs) {
    $rng = $para.Range
    if ($rng.Text -eq "`r" -and $rng.Style.NameLocal -eq "List Bullet" -and $rng.Font.Bold -eq -1 -and $rng.ParagraphFormat.LeftIndent -eq 0) {
        $prev = $para.Previous(1)
        if ($prev -ne $null -and $prev.Range.Text -like "*Space Complexity*") {
            $target = $para
            break
        }
    }
}

if ($target -eq $null) {
    Write-Output "Target paragraph not found"
} else {
    $r = $target.Range

    # Add the hanging-indent formatting: w:ind w:left="360" w:hanging="360"
    # (360 twips = 18 points; "hanging" is a negative first-line indent).
    $r.ParagraphFormat.LeftIndent = 18
    $r.ParagraphFormat.FirstLineIndent = -18

    # Turn off the complex-script bold flag (w:bCs) on the paragraph mark.
    $r.Font.BoldBi = 0

    # Turn off the regular bold flag (w:b) on the paragraph mark. Setting
    # Font.Bold directly on a completely run-less range isn't accepted
    # here, so briefly insert a placeholder character, flip Bold off
    # (which also updates the stored paragraph-mark run properties), then
    # delete the placeholder again, restoring the empty paragraph.
    $target.Range.Select()
    $sel = $word.Selection
    $sel.Collapse(1)
    $sel.TypeText("X")

    $afterType = $target.Range
    $afterType.Font.Bold = 0

    $placeholder = $d.Range($afterType.Start, $afterType.Start + 1)
    $placeholder.Delete()

    Write-Output "Updated paragraph. Text=[$($target.Range.Text)] Bold=$($target.Range.Font.Bold) LeftIndent=$($target.Range.ParagraphFormat.LeftIndent) FirstLineIndent=$($target.Range.ParagraphFormat.FirstLineIndent)"
}
